$d = $word.ActiveDocument

# Each "<id>p104r_N</id>" was originally split across three runs:
#   run1 (Courier New) "<id>"
#   run2 (Arial)        "p104r_N"
#   run3 (Courier New)  "</id>"
# Re-typing the same visible text via Find/Replace merges it back into a
# single run that inherits the formatting of the first matched run
# (Courier New / 7f6000 / sz 18), eliminating the stray Arial run in the
# middle, exactly as in the target document.
for ($i = 1; $i -le 5; $i++) {
    $text = "<id>p104r_$i</id>"
    $d.Content.Find.Execute($text, $true, $false, $false, $false, $false,
                             $true, 1, $false, $text, 2)
}
